# Add all columns from corpus
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
# D1 used to hold "Pays"; it becomes "Profession" and "Pays" moves to E1.
$ws.Range("D1").Value = "Profession"
$ws.Range("E1").Value = "Pays"
$ws.Range("F1").Value = "Titre"
$ws.Range("G1").Value = "Source"
$ws.Range("H1").Value = "Medium"
$ws.Range("I1").Value = "Publication"
$ws.Range("J1").Value = "Thème(s) critique(s) (idées primaires)"
$ws.Range("K1").Value = "Paradigme esthétique (idées primaires)"
$ws.Range("L1").Value = "Motifs esthétiques"
$ws.Range("M1").Value = "Registre(s)"
$ws.Range("N1").Value = "Interprétation"

# Header row gets taller to accommodate the new, longer titles.
$ws.Rows.Item(1).RowHeight = 94.5

# --- Row 2 (record 1) ---
# D2 used to hold "Australie"; it becomes "Blabla" and "Australie" moves to E2.
$ws.Range("D2").Value = "Blabla"
$ws.Range("E2").Value = "Australie"
$ws.Range("F2").Value = "Blabla"
$ws.Range("G2").Value = "Blabla"
$ws.Range("H2").Value = "Blabla"
$ws.Range("I2").Value = "Blabla"
$ws.Range("J2").Value = "Blabla"
$ws.Range("K2").Value = "Blabla"
$ws.Range("L2").Value = "Blabla"
$ws.Range("M2").Value = "Blabla"
$ws.Range("N2").Value = "Blabla"

# --- Row 3 (record 2) ---
# D3 used to hold "France"; it becomes "Boublou" and "France" moves to E3.
$ws.Range("D3").Value = "Boublou"
$ws.Range("E3").Value = "France"
$ws.Range("F3").Value = "Boublou"
$ws.Range("G3").Value = "Boublou"
$ws.Range("H3").Value = "Boublou"
$ws.Range("I3").Value = "Boublou"
$ws.Range("J3").Value = "Boublou"
$ws.Range("K3").Value = "Boublou"
$ws.Range("L3").Value = "Boublou"
$ws.Range("M3").Value = "Boublou"
$ws.Range("N3").Value = "Boublou"

# Match the final selection recorded in the workbook view.
[void]$ws.Range("N2:N3").Select()
